$d = $word.ActiveDocument

# The flyer's source list drops the "chinese-lantern" credit line. In the
# original document that credit lives in its own paragraph (a hyperlink run)
# immediately followed by a blank spacer paragraph, mirroring the pattern
# used between every other credit in the list. Removing both paragraphs
# keeps the remaining entries' spacing/pattern intact.
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "*chinese-lantern-3232431*") {
        $startPos = $p.Range.Start
        $endPos = $p.Range.End
        $nextPara = $p.Next()
        if ($nextPara -ne $null -and $nextPara.Range.Text -eq "`r") {
            $endPos = $nextPara.Range.End
        }
        $d.Range($startPos, $endPos).Delete()
        break
    }
}
